# Part_List_RISM.xlsx update
# Adds a new "Carte contacteur" parts table (rows 58-73) to Feuil1,
# mirroring the layout of the existing "Carte pédalier" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Section title (row 58) ---
$ws.Range("B58").Value = "Carte contacteur"
$ws.Range("B58").Font.Bold = $true

# --- Table header (row 59) ---
$ws.Range("A59").Value = "Type"
$ws.Range("B59").Value = "Valeur"
$ws.Range("C59").Value = "Référence EPSA"
$ws.Range("D59").Value = "Référence Farnell"
$ws.Range("E59").Value = "Emplacement"
$ws.Range("F59").Value = "Nb de pièces"
$ws.Range("G59").Value = "UdV"
$ws.Range("H59").Value = "Prix UdV"
$ws.Range("I59").Value = "Prix pour la carte"

# --- Data rows (60-73) ---
$ws.Range("A60").Value = "Condensateur"
$ws.Range("B60").Value = "10u"
$ws.Range("E60").Value = "C1"
$ws.Range("F60").Value = 1

$ws.Range("A61").Value = "Condensateur"
$ws.Range("B61").Value = "27p"
$ws.Range("E61").Value = "C2,C4"
$ws.Range("F61").Value = 2

$ws.Range("A62").Value = "Condensateur"
$ws.Range("B62").Value = "1u"
$ws.Range("E62").Value = "C3"
$ws.Range("F62").Value = 1

$ws.Range("A63").Value = "Diode"
$ws.Range("B63").Value = "DIODE ?R3"
$ws.Range("E63").Value = "D1,D2,D3"
$ws.Range("F63").Value = 3

$ws.Range("A64").Value = "Relais"
$ws.Range("B64").Value = "RELAIS-DPDT"
$ws.Range("E64").Value = "K1,K2,K3"
$ws.Range("F64").Value = 3

$ws.Range("A65").Value = "Connecteur"
$ws.Range("B65").Value = "KK2"
$ws.Range("E65").Value = "P1"
$ws.Range("F65").Value = 1

$ws.Range("A66").Value = "Connecteur"
$ws.Range("B66").Value = "KK4"
$ws.Range("E66").Value = "P2"
$ws.Range("F66").Value = 1

$ws.Range("A67").Value = "Connecteur"
$ws.Range("B67").Value = "Bornier 4"
$ws.Range("E67").Value = "P3"
$ws.Range("F67").Value = 1

$ws.Range("A68").Value = "Transistor"
$ws.Range("B68").Value = "FET_N"
$ws.Range("E68").Value = "Q1,Q2,Q3"
$ws.Range("F68").Value = 3

$ws.Range("A69").Value = "Résistance"
$ws.Range("B69").Value = "10k 0,25W"
$ws.Range("E69").Value = "R1,R2"
$ws.Range("F69").Value = 2

$ws.Range("A70").Value = "Résistance"
$ws.Range("B70").Value = "1k 0,25W"
$ws.Range("E70").Value = "R3,R5,R7"
$ws.Range("F70").Value = 3

$ws.Range("A71").Value = "Résistance"
$ws.Range("B71").Value = "220ohm 0,25W"
$ws.Range("E71").Value = "R4,R6,R8"
$ws.Range("F71").Value = 3

$ws.Range("A72").Value = "LM"
$ws.Range("B72").Value = "LM7805"
$ws.Range("E72").Value = "U1"
$ws.Range("F72").Value = 1

$ws.Range("A73").Value = "Logique"
$ws.Range("B73").Value = 4001
$ws.Range("B73").HorizontalAlignment = -4131
$ws.Range("E73").Value = "U5"
$ws.Range("F73").Value = 1

# --- Update the saved view/selection to match where the user ended up ---
[void]$ws.Range("F65").Select()
